$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: replace the MODE formula with the "Time Step" header text
$ws.Range("A1").Value = "Time Step"

# Fix a couple of values in the data block
$ws.Range("B3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("E12").Value = 0

# Add a new column H with per-row MODE formulas (set individually so
# each cell gets its own formula rather than a shared-formula group)
for ($r = 2; $r -le 12; $r++) {
    $ws.Range("H$r").Formula = "=MODE(C$r`:F$r)"
}

